$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 57. This shifts the existing row 57 (and all
# rows below it, down through the former row 158) down by one position,
# carrying their formatting along (so the date-formatted column D stays
# intact for every shifted row).
$ws.Rows("57:57").Insert()

# Populate the newly inserted row 57 with the new weekly price entry.
$ws.Cells.Item(57, 1).Value = 11
$ws.Cells.Item(57, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(57, 3).Value = "Bíobío"
$ws.Cells.Item(57, 4).Value = 44720
$ws.Cells.Item(57, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(57, 5).Value = 8
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100109
$ws.Cells.Item(57, 8).Value = "Uva"
$ws.Cells.Item(57, 9).Value = 100109001
$ws.Cells.Item(57, 10).Value = "Uva"
$ws.Cells.Item(57, 11).Value = "Red Globe"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 180
$ws.Cells.Item(57, 14).Value = 8000
$ws.Cells.Item(57, 15).Value = 9000
$ws.Cells.Item(57, 16).Value = 8556
$ws.Cells.Item(57, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(57, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(57, 19).Value = 856
$ws.Cells.Item(57, 20).Value = 10
